# contingencies with rene fine
# Extend the lines_states sheet with two more line contingencies (line7, line8).
# This shifts the existing extr1..extr8 rows down by two rows (now rows 10-17
# instead of 8-15), and refreshes the C/D/E (from_bus/to_bus/in_service)
# values for every row from 8 through 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the styled A-column format (border/bold/center, same as the other
# index cells) onto the two brand-new rows before filling their values.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(17, 1).PasteSpecial(-4122) # xlPasteFormats

# row, A (index), B (name), C (from_bus), D (to_bus), E (in_service)
$rows = @(
    @(2,  0,  "line1", 7,  9,  $true),
    @(3,  1,  "line2", 9,  8,  $false),
    @(4,  2,  "line3", 8,  10, $true),
    @(5,  3,  "line4", 8,  11, $true),
    @(6,  4,  "line5", 10, 5,  $true),
    @(7,  5,  "line6", 12, 8,  $true),
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
